$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.58"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("D3").Value = "'44.06"
$ws.Range("E3").Value = "'-0.45%"
$ws.Range("D4").Value = "'5.541"
$ws.Range("E4").Value = "'-0.56%"
$ws.Range("D5").Value = "'0.08116"
$ws.Range("E5").Value = "'0.24%"
$ws.Range("D6").Value = "'2.058"
$ws.Range("E6").Value = "'4.15%"
$ws.Range("D7").Value = "'0.9741"
$ws.Range("E7").Value = "'2.24%"
$ws.Range("D8").Value = "'0.1099"
$ws.Range("E8").Value = "'-6.51%"
$ws.Range("D9").Value = "'0.1894"
$ws.Range("E9").Value = "'2.12%"
$ws.Range("D10").Value = "'10.07"
$ws.Range("E10").Value = "'-1.62%"
$ws.Range("D11").Value = "'0.09967"
$ws.Range("E11").Value = "'0.04%"
$ws.Range("D12").Value = "'0.04727"
$ws.Range("E12").Value = "'-0.20%"
$ws.Range("E13").Value = "'-1.11%"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'-1.13%"
$ws.Range("D15").Value = "'0.04112"
$ws.Range("E15").Value = "'-2.61%"
$ws.Range("D16").Value = "'0.006090"
$ws.Range("E16").Value = "'2.78%"
$ws.Range("D17").Value = "'3.343"
$ws.Range("E17").Value = "'-0.84%"
$ws.Range("D18").Value = "'4.436"
$ws.Range("E18").Value = "'2.50%"
$ws.Range("D19").Value = "'2.612"
$ws.Range("E19").Value = "'1.66%"
$ws.Range("D20").Value = "'0.3345"
$ws.Range("E20").Value = "'-3.65%"
$ws.Range("D21").Value = "'0.1382"
$ws.Range("E21").Value = "'-2.01%"
$ws.Range("D22").Value = "'0.2572"
$ws.Range("E22").Value = "'2.65%"
$ws.Range("D23").Value = "'0.001304"
$ws.Range("E23").Value = "'4.45%"
$ws.Range("D24").Value = "'0.004384"
$ws.Range("E24").Value = "'1.10%"
$ws.Range("D25").Value = "'0.0001279"
$ws.Range("E25").Value = "'7.52%"
$ws.Range("D26").Value = "'0.0003737"
$ws.Range("E26").Value = "'-6.07%"
$ws.Range("E38").Value = "'0.60%"
$ws.Range("D39").Value = "'0.05635"
$ws.Range("E39").Value = "'1.64%"
$ws.Range("D40").Value = "'0.007597"
$ws.Range("E40").Value = "'0.86%"
$ws.Range("D41").Value = "'0.1415"
$ws.Range("E41").Value = "'0.36%"
$ws.Range("D42").Value = "'0.007545"
$ws.Range("E42").Value = "'-6.44%"
$ws.Range("D43").Value = "'0.001957"
$ws.Range("E43").Value = "'-2.87%"
$ws.Range("D44").Value = "'0.008315"
$ws.Range("E44").Value = "'-6.50%"
$ws.Range("D45").Value = "'0.00007054"
$ws.Range("E45").Value = "'-2.48%"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("D47").Value = "'0.0005797"
$ws.Range("E47").Value = "'-0.25%"
$ws.Range("D48").Value = "'0.002518"
$ws.Range("E48").Value = "'10.94%"
$ws.Range("D49").Value = "'0.003538"
$ws.Range("E49").Value = "'0.52%"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.06%"
